$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ends with a blank "template" row (170) below the last real
# entry (169), followed by three summary rows (sum [min], sum [h],
# sum [working weeks]). A new day of hours was logged, so:
#   1. The previous last entry (row 169) had its "end" time corrected.
#   2. The blank template row (170) is filled in with the new entry and
#      becomes a real data row, with F/G formulas extended into it.
#   3. A fresh blank template row is inserted below it (pushing the
#      summary rows down by one) so the sheet keeps its usual shape.
#   4. The "sum [min]" total is widened to include the new row.

# 1. Correct the end time recorded for row 169.
$ws.Range("E169").Value = 0.38541666666666669

# 2. Insert a new blank template row below row 170 FIRST (while row 170
#    still only carries its original, empty D/E/F formatting) so the
#    blank row's formatting mirrors the old template row exactly, then
#    push the summary rows (sum [min]/sum [h]/sum [working weeks]) down
#    from 171-173 to 172-174.
$ws.Rows("171").Insert()

# 3. Populate row 170 with the new day's data.
$ws.Range("A170").Value = 2014
$ws.Range("B170").Value = 7
$ws.Range("C170").Value = 31
$ws.Range("D170").Value = 0.59375
$ws.Range("E170").Value = 0.75
$ws.Range("F170").Formula = "=(E170-D170)*24*60"
$ws.Range("G170").Formula = "=F170/60"

# 4. Extend the "sum [min]" total to cover the new data row.
$ws.Range("F172").Formula = "=SUM(F2:F170)"

# Match the author's final active selection.
$ws.Range("F170").Select()
